$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 340 - this shifts the existing rows 340-369 down to
# 341-370 (matching the diff, which re-numbers every row from 340 onward
# by +1 and appends a brand-new row 370 at the end).
$ws.Rows.Item(340).Insert()

# Populate the newly inserted row 340 with the new weekly price record.
# (Columns A, B, C, E, F, G, H, I, J, K, R are constant for this whole
# sheet / this product's quality grades, same as every other row.)
$ws.Range("A340").Value2 = 10
$ws.Range("B340").Value2 = "Vega Modelo de Temuco"
$ws.Range("C340").Value2 = "La Araucanía"
$ws.Range("D340").Value2 = 44578
$ws.Range("E340").Value2 = 9
$ws.Range("F340").Value2 = "Fruta"
$ws.Range("G340").Value2 = 100108
$ws.Range("H340").Value2 = "Tropicales y subtropicales"
$ws.Range("I340").Value2 = 100108005
$ws.Range("J340").Value2 = "Piña"
$ws.Range("K340").Value2 = "Caramelo"
$ws.Range("L340").Value2 = "Primera"
$ws.Range("M340").Value2 = 75
$ws.Range("N340").Value2 = 19000
$ws.Range("O340").Value2 = 19000
$ws.Range("P340").Value2 = 19000
$ws.Range("Q340").Value2 = "$/caja 12 unidades"
$ws.Range("R340").Value2 = "Ecuador"
$ws.Range("S340").Value2 = 1583
$ws.Range("T340").Value2 = 12
